$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# Add new row of data: TUHH_MICROELECTRONICS with default validation value "Yes"
$ws.Range("A8").Value = "TUHH_MICROELECTRONICS"
$ws.Range("B8").Value = "Yes"

# Extend the data validation list that currently covers B1:B7 to also include B8
$ws.Range("B1:B8").Validation.Delete()
$ws.Range("B1:B8").Validation.Add(3, 1, 1, """Yes,No""")
$ws.Range("B1:B8").Validation.IgnoreBlank = $true
$ws.Range("B1:B8").Validation.InCellDropdown = $true
$ws.Range("B1:B8").Validation.ShowInput = $true
$ws.Range("B1:B8").Validation.ShowError = $true

# Update the active selection to B2, matching the post-edit state
$ws.Range("B2").Select()
